# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D occasionally look like plain numbers (e.g. '304.58'); the
# source workbook stores them as TEXT (general format, no numeric coercion).
# A leading apostrophe forces Excel/COM to keep them as text, matching the
# original inlineStr cell type instead of silently turning them into numbers.

$ws.Range('D2').Value = '43.329.24'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.319.95'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''304.58'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('D6').Value = '''100.97'
$ws.Range('E6').Value = '  -3.95%  '
$ws.Range('E7').Value = '  -3.39%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('D10').Value = '''35.26'
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D14').Value = '2.683.94'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '''15.60'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '2.337.28'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '43.273.63'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '''11.77'
$ws.Range('E19').Value = '  -1.33%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').Value = '''6.08'
$ws.Range('E21').Value = '  -2.75%  '
$ws.Range('D22').Value = '''68.19'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '''237.24'
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('E24').Value = '  -3.53%  '
$ws.Range('E25').Value = '  -3.81%  '
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').Value = '''24.71'
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('D28').Value = '''2.17'
$ws.Range('E28').Value = '  -1.93%  '
$ws.Range('D29').Value = '''34.53'
$ws.Range('E29').Value = '  -5.35%  '
$ws.Range('D30').Value = '''163.80'
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('E31').Value = '  -4.47%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '''2.42'
$ws.Range('E34').Value = '  -4.98%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '''4.54'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('D36').Value = '''16.90'
$ws.Range('E36').Value = '  -7.27%  '
$ws.Range('D37').Value = '''0.0704'
$ws.Range('E37').Value = '  -4.66%  '
$ws.Range('D38').Value = '''2.91'
$ws.Range('E38').Value = '  -4.74%  '
$ws.Range('E39').Value = '  -4.08%  '
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('D42').Value = '''2.57'
$ws.Range('E42').Value = '  +5.08%  '
$ws.Range('D43').Value = '1.971.15'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('D45').Value = '''18.49'
$ws.Range('E45').Value = '  -4.94%  '
$ws.Range('D46').Value = '''10.22'
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').Value = '''2.90'
$ws.Range('E47').Value = '  -5.54%  '
$ws.Range('D48').Value = '''55.78'
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '''1.55'
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.547.96'
$ws.Range('E51').Value = '  -0.40%  '
